$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A173").Value = "Appeals casework portal documentation"
$ws.Range("A173").WrapText = $true

$ws.Range("A174").Value = "Planning system"

$ws.Range("A173").Select() | Out-Null
